# Actualización automática 2025-05-29 08:40:07
#
# A new client row ("DISALME CIA. LTDA.") is inserted at row 10 on both
# worksheets ("VENTAS POR GRUPO" and "VENTA MENSUAL"), pushing every
# following client row down by one. The new row carries zero sales figures.
# The trailing summary row (counts / totals) then moves down one row too,
# and on the first sheet its "X de 25" labels become "X de 26" because the
# total number of clients grew from 25 to 26.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Sheet 1: "VENTAS POR GRUPO"  (columns A:N, data rows 2-26, summary row 27)
# ---------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("VENTAS POR GRUPO")

# Insert a whole new row before row 10; this shifts rows 10-27 down to 11-28
# and keeps styles/dimension in sync automatically.
$ws1.Rows.Item(10).Insert()

# Fill in the new client row (row 10) with the asesor name, the new
# client name, and zero for every product column.
$ws1.Cells.Item(10, 1).Value = "LOZANO MOLINA TITO"
$ws1.Cells.Item(10, 2).Value = "DISALME CIA. LTDA."
for ($col = 3; $col -le 14; $col++) {
    $ws1.Cells.Item(10, $col).Value = 0
}

# The summary row (now row 28) previously read "... de 25" - there are now
# 26 clients, so update each label accordingly while keeping the counts.
$ws1.Cells.Item(28, 3).Value  = "0 de 26"
$ws1.Cells.Item(28, 4).Value  = "1 de 26"
$ws1.Cells.Item(28, 5).Value  = "0 de 26"
$ws1.Cells.Item(28, 6).Value  = "0 de 26"
$ws1.Cells.Item(28, 7).Value  = "0 de 26"
$ws1.Cells.Item(28, 8).Value  = "0 de 26"
$ws1.Cells.Item(28, 9).Value  = "0 de 26"
$ws1.Cells.Item(28, 10).Value = "0 de 26"
$ws1.Cells.Item(28, 11).Value = "1 de 26"
$ws1.Cells.Item(28, 12).Value = "2 de 26"
$ws1.Cells.Item(28, 13).Value = "0 de 26"
$ws1.Cells.Item(28, 14).Value = "0 de 26"

# ---------------------------------------------------------------------
# Sheet 2: "VENTA MENSUAL"  (columns A:F, data rows 2-26, summary row 27)
# ---------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("VENTA MENSUAL")

# Same row insertion as sheet 1.
$ws2.Rows.Item(10).Insert()

$ws2.Cells.Item(10, 1).Value = "LOZANO MOLINA TITO"
$ws2.Cells.Item(10, 2).Value = "DISALME CIA. LTDA."
for ($col = 3; $col -le 6; $col++) {
    $ws2.Cells.Item(10, $col).Value = 0
}

# The totals row (now row 28) keeps the same sums - the inserted row only
# contributed zeros - so nothing else needs to change there.
